$d = $word.ActiveDocument

# 1) Merge "...calls and " + "that" + " " into a single run ending in "...calls and that "
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Replacement.ClearFormatting()
$r1.Find.Execute("calls and that ", $true, $false, $false, $false, $false, $true, 1, $false, "calls and that ", 2) | Out-Null

# 2) Replace the two runs "6" and "2" (rendered as "62") with a single run "55"
#    Insert "55" right at the boundary between the "6" run and the "2" run so it
#    inherits the (matching) formatting of those runs, then delete the old "6"/"2".
$r2 = $d.Content
$r2.Find.ClearFormatting()
$found2 = $r2.Find.Execute("62", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $start2 = $r2.Start
    $mid2 = $d.Range($start2 + 1, $start2 + 1)
    $mid2.InsertBefore("55")
    $d.Range($start2, $start2 + 1).Delete() | Out-Null
    $d.Range($start2 + 2, $start2 + 3).Delete() | Out-Null
}

# 3) Merge "...required for " + "each" + " " into a single run ending in "...required for each "
$r3 = $d.Content
$r3.Find.ClearFormatting()
$r3.Find.Replacement.ClearFormatting()
$r3.Find.Execute("required for each ", $true, $false, $false, $false, $false, $true, 1, $false, "required for each ", 2) | Out-Null

# 4) Merge "...traffic intensity " + "is" + " doubled?" into a single run
$r4 = $d.Content
$r4.Find.ClearFormatting()
$r4.Find.Replacement.ClearFormatting()
$r4.Find.Execute("traffic intensity is doubled?", $true, $false, $false, $false, $false, $true, 1, $false, "traffic intensity is doubled?", 2) | Out-Null
